# Rename "FP"/"TP" sheets to "NL"/"LL" (NLRating/LLRating terminology),
# update their rating-column headers, drop the now-unused extra columns
# (ReaderID/ModalityID/Paradigm) from the TRUTH sheet, and restore the
# per-sheet active-cell selections to match the saved file.

$wb = $excel.ActiveWorkbook

# --- "FP" -> "NL" --------------------------------------------------------
$wsNL = $wb.Worksheets.Item("FP")
$wsNL.Name = "NL"
$wsNL.Range("D1").Value = "NLRating"
[void]$wsNL.Range("D2").Select()

# --- "TP" -> "LL" --------------------------------------------------------
$wsLL = $wb.Worksheets.Item("TP")
$wsLL.Name = "LL"
$wsLL.Range("E1").Value = "LLRating"
[void]$wsLL.Range("G7").Select()

# --- TRUTH: drop ReaderID/ModalityID/Paradigm columns -------------------
$wsTruth = $wb.Worksheets.Item("TRUTH")
$wsTruth.Range("D1:F13").Delete()
[void]$wsTruth.Range("C2").Select()
